# Applies the "Penalty Reward System" (unfinished) edit:
#  - Sheet "Forecast Comparison": shift Week_Start_Date values forward by one
#    week and overwrite MyForecast (column D) with new values.
#  - Sheet "Summary": update Historical Range, Total Forecast / Max-Min
#    Forecast figures to reflect the new MyForecast column.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (column B) for rows 2..17.
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (column D) for rows 2..17.
$newForecast = @(0, 0, 0, 0, 0, 0, 1, 1, 1, 1, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $cellB = $wsForecast.Range("B$row")
    $cellB.NumberFormat = "@"
    $cellB.Value = $newDates[$i]

    $wsForecast.Range("D$row").Value = $newForecast[$i]
}

# Summary sheet updates.
$cellHist = $wsSummary.Range("B2")
$cellHist.NumberFormat = "@"
$cellHist.Value = "2023-01-08 to 2025-01-05"

$cellB9 = $wsSummary.Range("B9")
$cellB9.NumberFormat = "@"
$cellB9.Value = "6"

$cellB10 = $wsSummary.Range("B10")
$cellB10.NumberFormat = "@"
$cellB10.Value = "4"

$cellB11 = $wsSummary.Range("B11")
$cellB11.NumberFormat = "@"
$cellB11.Value = "2"

$cellB12 = $wsSummary.Range("B12")
$cellB12.NumberFormat = "@"
$cellB12.Value = "1"

$cellB14 = $wsSummary.Range("B14")
$cellB14.NumberFormat = "@"
$cellB14.Value = "0"

$cellB15 = $wsSummary.Range("B15")
$cellB15.NumberFormat = "@"
$cellB15.Value = "2025-03-30"
